$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2088.25
$ws.Range("I100").Value = 1920
$ws.Range("J100").Value = 2368.6667
$ws.Range("K100").Value = 1920
$ws.Range("L100").Value = 2368.6667
$ws.Range("M100").Value = -1379
$ws.Range("N100").Value = -3450.6667

$ws.Range("H127").Value = 1622.9412
$ws.Range("I127").Value = 719.4
$ws.Range("J127").Value = 1999.4166
$ws.Range("K127").Value = 2158.2
$ws.Range("L127").Value = 5998.2498
$ws.Range("M127").Value = 2801.8
$ws.Range("N127").Value = -15918.2498

$ws.Range("H128").Value = 40000
$ws.Range("I128").Value = 40000
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 40000
$ws.Range("L128").Value = 0
$ws.Range("M128").Value = -35020
$ws.Range("N128").ClearContents()

$ws.Range("H129").Value = 773.7143
$ws.Range("I129").Value = 559.6667
$ws.Range("J129").Value = 1159
$ws.Range("K129").Value = 1679.0001
$ws.Range("L129").Value = 3477
$ws.Range("M129").Value = 3320.9999
$ws.Range("N129").Value = -13477

$ws.Range("H132").Value = 682281.5
$ws.Range("I132").Value = 1523.0834
$ws.Range("J132").Value = 4086073.8
$ws.Range("K132").Value = 4569.2502
$ws.Range("L132").Value = 12258221.4
$ws.Range("M132").Value = -2039.2502
$ws.Range("N132").Value = -12263281.4

$ws.Range("H135").Value = 22436.709
$ws.Range("I135").Value = 28992.945
$ws.Range("J135").Value = 2768
$ws.Range("K135").Value = 260936.505
$ws.Range("L135").Value = 24912
$ws.Range("M135").Value = -258401.505
$ws.Range("N135").Value = -29982

$ws.Range("H137").Value = 2001361.4
$ws.Range("I137").Value = 3226922.2
$ws.Range("J137").Value = 1762.1052
$ws.Range("K137").Value = 9680766.600000001
$ws.Range("L137").Value = 5286.3156
$ws.Range("M137").Value = -9678216.600000001
$ws.Range("N137").Value = -10386.3156

$ws.Range("H138").Value = 2034347.5
$ws.Range("I138").Value = 1255.3658
$ws.Range("J138").Value = 4067439.8
$ws.Range("K138").Value = 3766.0974
$ws.Range("L138").Value = 12202319.4
$ws.Range("M138").Value = 1373.9026
$ws.Range("N138").Value = -12212599.4

$ws.Range("H141").Value = 1283.1041
$ws.Range("I141").Value = 1308.7273
$ws.Range("J141").Value = 1001.25
$ws.Range("K141").Value = 3926.1819
$ws.Range("L141").Value = 3003.75
$ws.Range("M141").Value = 1253.8181
$ws.Range("N141").Value = -13363.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 3103.5
$ws.Range("I14").Value = 1200
$ws.Range("J14").Value = 5007
$ws.Range("K14").Value = 1200
$ws.Range("L14").Value = 5007
$ws.Range("M14").Value = -1025
$ws.Range("N14").Value = -5357

$ws.Range("H32").Value = 15981.506
$ws.Range("I32").Value = 20861.322
$ws.Range("J32").Value = 6384.533
$ws.Range("K32").Value = 20861.322
$ws.Range("L32").Value = 6384.533
$ws.Range("M32").Value = -20574.322
$ws.Range("N32").Value = -6958.533

$ws.Range("H61").Value = 19647908
$ws.Range("I61").Value = 22245328
$ws.Range("J61").Value = 167250
$ws.Range("K61").Value = 22245328
$ws.Range("L61").Value = 167250
$ws.Range("M61").Value = -22245116
$ws.Range("N61").Value = -167674

$ws.Range("H97").Value = 2841820
$ws.Range("I97").Value = 3907083.8
$ws.Range("J97").Value = 1116.6666
$ws.Range("K97").Value = 3907083.8
$ws.Range("L97").Value = 1116.6666
$ws.Range("M97").Value = -3906587.8
$ws.Range("N97").Value = -2108.6666

$ws.Range("H122").Value = 22225140
$ws.Range("I122").Value = 3648
$ws.Range("J122").Value = 111111110
$ws.Range("K122").Value = 10944
$ws.Range("L122").Value = 333333330
$ws.Range("M122").Value = -8494
$ws.Range("N122").Value = -333338230

$ws.Range("H136").Value = 19647908
$ws.Range("I136").Value = 22245328
$ws.Range("J136").Value = 167250
$ws.Range("K136").Value = 66735984
$ws.Range("L136").Value = 501750
$ws.Range("M136").Value = -66733434
$ws.Range("N136").Value = -506850

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2357
$ws.Range("I134").Value = 1265.85
$ws.Range("J134").Value = 4175.5835
$ws.Range("K134").Value = 3797.55
$ws.Range("L134").Value = 12526.7505
$ws.Range("M134").Value = -1262.55
$ws.Range("N134").Value = -17596.7505

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2517.5098
$ws.Range("I31").Value = 1019.8
$ws.Range("J31").Value = 5793.75
$ws.Range("K31").Value = 1019.8
$ws.Range("L31").Value = 5793.75
$ws.Range("M31").Value = -724.8
$ws.Range("N31").Value = -6383.75

$ws.Range("H34").Value = 2517.5098
$ws.Range("I34").Value = 1019.8
$ws.Range("J34").Value = 5793.75
$ws.Range("K34").Value = 1019.8
$ws.Range("L34").Value = 5793.75
$ws.Range("M34").Value = -817.8
$ws.Range("N34").Value = -6197.75

$ws.Range("H99").Value = 8559.200000000001
$ws.Range("I99").Value = 9014.154
$ws.Range("J99").Value = 7714.2856
$ws.Range("K99").Value = 9014.154
$ws.Range("L99").Value = 7714.2856
$ws.Range("M99").Value = -7516.154
$ws.Range("N99").Value = -10710.2856

$ws.Range("H126").Value = 8559.200000000001
$ws.Range("I126").Value = 9014.154
$ws.Range("J126").Value = 7714.2856
$ws.Range("K126").Value = 27042.462
$ws.Range("L126").Value = 23142.8568
$ws.Range("M126").Value = -24572.462
$ws.Range("N126").Value = -28082.8568

$ws.Range("H132").Value = 19665.76
$ws.Range("I132").Value = 22890.195
$ws.Range("J132").Value = 1125.25
$ws.Range("K132").Value = 68670.58499999999
$ws.Range("L132").Value = 3375.75
$ws.Range("M132").Value = -66140.58499999999
$ws.Range("N132").Value = -8435.75

$ws.Range("H134").Value = 19907.965
$ws.Range("I134").Value = 1230.9546
$ws.Range("J134").Value = 78607.14
$ws.Range("K134").Value = 3692.8638
$ws.Range("L134").Value = 235821.42
$ws.Range("M134").Value = -1157.8638
$ws.Range("N134").Value = -240891.42

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 221.09091
$ws.Range("I2").Value = 248
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 1488
$ws.Range("L2").Value = 600
$ws.Range("M2").Value = -1375
$ws.Range("N2").Value = -826

$ws.Range("H15").Value = 547.125
$ws.Range("I15").Value = 110
$ws.Range("J15").Value = 648
$ws.Range("K15").Value = 330
$ws.Range("L15").Value = 1944
$ws.Range("M15").Value = -190
$ws.Range("N15").Value = -2224

$ws.Range("H131").Value = 1050.7826
$ws.Range("I131").Value = 406.9091
$ws.Range("J131").Value = 1253.1428
$ws.Range("K131").Value = 1220.7273
$ws.Range("L131").Value = 3759.4284
$ws.Range("M131").Value = 3819.2727
$ws.Range("N131").Value = -13839.4284

$ws.Range("H132").Value = 1871.0714
$ws.Range("I132").Value = 1711.25
$ws.Range("J132").Value = 2084.1667
$ws.Range("K132").Value = 15401.25
$ws.Range("L132").Value = 18757.5003
$ws.Range("M132").Value = -12871.25
$ws.Range("N132").Value = -23817.5003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5228.8335
$ws.Range("I122").Value = 4935.6665
$ws.Range("J122").Value = 5522
$ws.Range("K122").Value = 14806.9995
$ws.Range("L122").Value = 16566
$ws.Range("M122").Value = -12356.9995
$ws.Range("N122").Value = -21466

$ws.Range("H135").Value = 30000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 30000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 30000
$ws.Range("N135").Value = -40140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3447.077
$ws.Range("I40").Value = 3223
$ws.Range("J40").Value = 3951.25
$ws.Range("K40").Value = 3223
$ws.Range("L40").Value = 3951.25
$ws.Range("M40").Value = -3087
$ws.Range("N40").Value = -4223.25

$ws.Range("H132").Value = 30918.7
$ws.Range("I132").Value = 20913.473
$ws.Range("J132").Value = 62111.47
$ws.Range("K132").Value = 62740.41900000001
$ws.Range("L132").Value = 186334.41
$ws.Range("M132").Value = -60210.41900000001
$ws.Range("N132").Value = -191394.41

$ws.Range("H136").Value = 127782.82
$ws.Range("I136").Value = 101341.63
$ws.Range("J136").Value = 176258.33
$ws.Range("K136").Value = 304024.89
$ws.Range("L136").Value = 528774.99
$ws.Range("M136").Value = -301474.89
$ws.Range("N136").Value = -533874.99

$ws.Range("H137").Value = 25162.533
$ws.Range("I137").Value = 20390
$ws.Range("J137").Value = 25896.77
$ws.Range("K137").Value = 20390
$ws.Range("L137").Value = 25896.77
$ws.Range("M137").Value = -15290
$ws.Range("N137").Value = -36096.77

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 54217.742
$ws.Range("I132").Value = 46196.59
$ws.Range("J132").Value = 79427.07000000001
$ws.Range("K132").Value = 138589.77
$ws.Range("L132").Value = 238281.21
$ws.Range("M132").Value = -136059.77
$ws.Range("N132").Value = -243341.21

$ws.Range("H136").Value = 136213.27
$ws.Range("I136").Value = 501000
$ws.Range("J136").Value = 80092.234
$ws.Range("K136").Value = 1503000
$ws.Range("L136").Value = 240276.702
$ws.Range("M136").Value = -1500450
$ws.Range("N136").Value = -245376.702
